$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Columns A-D (1-4) hold image filenames like "face//face_26.jpg" that
# should become "book//book_26.jpg". Column L (12) holds single-letter
# position codes (b/y/r) that should be expanded to full words.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string]) {
            if ($val.Contains("face")) {
                $cell.Value = $val.Replace("face", "book")
            }
        }
    }

    $lCell = $ws.Cells.Item($r, 12)
    $lVal = $lCell.Value2
    if ($lVal -is [string]) {
        if ($lVal -eq "b") {
            $lCell.Value = "center"
        } elseif ($lVal -eq "y") {
            $lCell.Value = "left"
        } elseif ($lVal -eq "r") {
            $lCell.Value = "right"
        }
    }
}
